# Scenario3_savepath/Assets_CPY_input.xlsx
# Add a new reservoir row ("Kaeng Suea Ten") as row 7 of the data table,
# matching the pattern of the existing rows (shared formula for column D).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "Kaeng Suea Ten"
$ws.Range("C7").Value = 48
$ws.Range("D7").Formula = "=C7*1000000/E7/1000/3600*86400*30.25/1000000"
$ws.Range("E7").Value = 0.177
$ws.Range("F7").Value = 1175
